$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Overwrite row 3 ("S1") with the values that used to live in row 4 ("Test")
$ws.Range("A3").Value = "Test"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 0.1
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 3

# Remove the now-duplicate row 4, shrinking the used range back to A1:J3
$ws.Rows.Item(4).Delete()
